# Deploy update from 15e2cca904af7fba09cb3c8ac702eecd65e20b79
#
# Converts the eight numbered "Heading 3" section titles into plain
# Body Text paragraphs prefixed with "###" (Markdown-style heading
# markers), removes the per-section bookmarks that used to anchor
# those headings, promotes each section's "First Paragraph" body copy
# to plain "Body Text", and re-seats the top-level document bookmark
# so it wraps the title run (bookmarkStart before the run, bookmarkEnd
# right after it) instead of being a zero-width bookmark ahead of it.

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Re-seat the "data-management-plan-template" bookmark so it
#    spans the whole "Data Management Plan Template" title run
#    (bookmarkStart ... run ... bookmarkEnd) rather than sitting as
#    a zero-width bookmark ahead of the run.
# ---------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range.Duplicate
# drop the trailing paragraph mark so the bookmark hugs just the text
$titleRange.MoveEndWhile(" `r`n", -1) | Out-Null
$titleRange.End = $titleRange.End - 1
$d.Bookmarks("data-management-plan-template").Delete()
$d.Bookmarks.Add("data-management-plan-template", $titleRange) | Out-Null

# ---------------------------------------------------------------
# 2. Turn each "N. <Section Title>" Heading 3 paragraph into a Body
#    Text paragraph with a leading "###" marker, drop the bookmark
#    that used to mark that section heading, and promote the body
#    copy paragraph right after it from "First Paragraph" to plain
#    "Body Text". Headings are located via their bookmarks so the
#    script does not depend on hard-coded paragraph indices.
# ---------------------------------------------------------------
$sectionBookmarks = @(
    "data-formats",
    "metadata",
    "documentation",
    "policies-for-access-and-sharing",
    "policies-and-provisions-for-reuse-and-redistribution",
    "plans-for-archiving-and-preservation",
    "versioning-of-stored-assets",
    "data-security"
)

foreach ($bookmarkName in $sectionBookmarks) {
    $bookmark = $d.Bookmarks($bookmarkName)
    $headingPara = $bookmark.Range.Paragraphs(1)
    $bodyPara = $headingPara.Next()

    $bookmark.Delete()

    $headingPara.Range.InsertBefore("###")
    $headingPara.Style = "Body Text"

    $bodyPara.Style = "Body Text"
}

Write-Output "done"
